# 2023 - Day 16 - Solved part II
# Fill in the Part I / Part II timings (ms) for Day 16 (row 21).
# Formulas in column G (Part I & II) and the Average/Total summary rows
# (32/34) recalculate automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D21").Value = 8
$ws.Range("E21").Value = 1973

# Match the author's on-screen selection after entering the Day 16 data.
$ws.Range("D22").Select()
